$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.805.11"
$ws.Range("E2").Value = "  +1.95%  "

$ws.Range("D3").Value = "3.003.63"
$ws.Range("E3").Value = "  +1.07%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "511.63"
$ws.Range("E5").Value = "  +4.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.35"
$ws.Range("E6").Value = "  +4.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +3.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.46"
$ws.Range("E9").Value = "  +4.70%  "

$ws.Range("E10").Value = "  +6.47%  "

$ws.Range("E11").Value = "  +1.79%  "

$ws.Range("E12").Value = "  +1.88%  "

$ws.Range("D13").Value = "3.518.42"
$ws.Range("E13").Value = "  +1.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.99"
$ws.Range("E14").Value = "  +4.01%  "

$ws.Range("E15").Value = "  +10.48%  "

$ws.Range("D16").Value = "56.835.84"
$ws.Range("E16").Value = "  +2.33%  "

$ws.Range("D17").Value = "3.003.24"
$ws.Range("E17").Value = "  +1.53%  "

$ws.Range("E18").Value = "  +4.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.55"
$ws.Range("E19").Value = "  +2.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.85"
$ws.Range("E20").Value = "  +4.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.97"
$ws.Range("E21").Value = "  +1.26%  "

$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("E23").Value = "  +4.70%  "

$ws.Range("E24").Value = "  +4.91%  "

$ws.Range("E25").Value = "  +4.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("E27").Value = "  +7.04%  "

$ws.Range("E28").Value = "  +0.80%  "

$ws.Range("E29").Value = "  +6.98%  "

$ws.Range("E30").Value = "  +3.68%  "

$ws.Range("E31").Value = "  +6.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.57"
$ws.Range("E32").Value = "  +5.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.97"
$ws.Range("E33").Value = "  +4.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.58"
$ws.Range("E34").Value = "  +3.58%  "

$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("E36").Value = "  -2.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0680"
$ws.Range("E37").Value = "  +4.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.96"
$ws.Range("E38").Value = "  +1.80%  "

$ws.Range("D39").Value = "3.035.95"
$ws.Range("E39").Value = "  +1.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.11"
$ws.Range("E40").Value = "  +2.20%  "

$ws.Range("E41").Value = "  +0.14%  "

$ws.Range("D42").Value = "2.294.51"
$ws.Range("E42").Value = "  +7.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.648"
$ws.Range("E43").Value = "  +2.24%  "

$ws.Range("E44").Value = "  +3.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  -0.24%  "

$ws.Range("E46").Value = "  +2.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.95"
$ws.Range("E47").Value = "  +9.17%  "

$ws.Range("E48").Value = "  +5.04%  "

$ws.Range("E49").Value = "  +1.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.40"
$ws.Range("E50").Value = "  -0.18%  "

$ws.Range("E51").Value = "  +4.38%  "
